# ---------------------------------------------------------------------------
# Deck update: "Updated deck added some comments"
#  1. Bump the auto date-placeholder text (10/16/16 -> 10/21/16) everywhere
#     it appears: slide master, every slide layout, the handout master and
#     the notes master.
#  2. Slide 3  : tidy the "npm gulp" / "npm run mocha" demo commands.
#  3. Slide 5  : tidy the title + add detail to the "Bootstraps" bullet.
#  4. Slide 13 : merge the "Use gulp-jspm task to perform " run fragments.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Date placeholder text -----------------------------------------------
$oldDate = "10/16/16"
$newDate = "10/21/16"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master itself
Update-DateShape $p.SlideMaster.Shapes

# Every layout hanging off the slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShape $layouts.Item($li).Shapes
}

# Handout master
Update-DateShape $p.HandoutMaster.Shapes

# Notes master
Update-DateShape $p.NotesMaster.Shapes

# --- 2. Slide 3 : demo command bullets --------------------------------------
$slide3 = $p.Slides.Item(3)
$content3 = $slide3.Shapes.Item(3)
$tr3 = $content3.TextFrame.TextRange

# "npm " + "gulp" -> "gulp"
$para2 = $tr3.Paragraphs(2, 1)
$para2.Text = "gulp"

# "npm run mocha" -> "npm run " + "coverage"
$para3 = $tr3.Paragraphs(3, 1)
$para3.Text = "npm run "
$para3.InsertAfter("coverage") | Out-Null

# --- 3. Slide 5 : title + bootstraps bullet ---------------------------------
$slide5 = $p.Slides.Item(5)

$title5 = $slide5.Shapes.Item(1)
$trTitle = $title5.TextFrame.TextRange
$paraTitle = $trTitle.Paragraphs(1, 1)
$paraTitle.Text = "High Level Choices (tmp)"
$paraTitle = $trTitle.Paragraphs(1, 1)
$paraTitle.Text = "High Level Choices"

$content5 = $slide5.Shapes.Item(2)
$trContent = $content5.TextFrame.TextRange
$paraBoot = $trContent.Paragraphs(5, 1)
$paraBoot.Text = "Bootstraps "
$paraBoot.InsertAfter("differently than Angular 1.x") | Out-Null

# --- 4. Slide 13 : merge the gulp-jspm task sentence fragments --------------
$slide13 = $p.Slides.Item(13)
$content13 = $slide13.Shapes.Item(2)
$tr13 = $content13.TextFrame.TextRange
$para5_13 = $tr13.Paragraphs(5, 1)

$mergedText = "Use gulp-jspm task to perform "
$sub = $para5_13.Characters(1, $mergedText.Length)
$sub.Text = "Use gulp-jspm task to perform (tmp)"
$para5_13 = $tr13.Paragraphs(5, 1)
$sub = $para5_13.Characters(1, "Use gulp-jspm task to perform (tmp)".Length)
$sub.Text = $mergedText
